{"js": "// Load all paragraphs in the document body so we can locate the\n// \"Baz changes\" paragraph (which also carries the \"_GoBack\" bookmark)\n// and the trailing empty paragraphs.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Find the paragraph whose text is \"Baz changes\" (it is split across\n// two runs with a bookmark in between, but Paragraph.text concatenates\n// the runs so it still reads \"Baz changes\").\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Baz changes\") {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!targetParagraph) {\n  throw new Error('Could not find the \"Baz changes\" paragraph.');\n}\n\n// The \"_GoBack\" bookmark currently sits inside that paragraph; remove\n// it from there so it can be re-created at the end of the document.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Replace the paragraph's text with the new sentence and make it bold\n// and red, matching both the run formatting and the paragraph mark\n// formatting.\ntargetParagraph.clear();\ntargetParagraph.font.bold = true;\ntargetParagraph.font.color = \"#FF0000\";\ntargetParagraph.insertText(\n  \"Version management help many developers to collaborate while working on same project.\",\n  Word.InsertLocation.start\n);\n\n// Re-insert the \"_GoBack\" bookmark at the last paragraph in the body.\nconst lastParagraph = body.paragraphs.getLast();\nlastParagraph.getRange().insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Replace the \"Baz changes\" paragraph with a new bold, red sentence\n# about version management, and relocate the \"_GoBack\" bookmark that\n# used to sit inside that paragraph to the final (empty) paragraph of\n# the document.\n$d = $word.ActiveDocument\n\n# The \"Baz changes\" paragraph is the 5th paragraph in the body.\n$targetParagraph = $d.Paragraphs(5)\n$targetRange = $targetParagraph.Range\n\n# Exclude the trailing paragraph mark from the range we overwrite so we\n# only replace the visible text (this also removes the \"_GoBack\"\n# bookmark that lived inside this paragraph, since its anchor text is\n# being deleted).\n$wdCharacter = [Microsoft.Office.Interop.Word.WdUnits]::wdCharacter\n$targetRange.MoveEnd($wdCharacter, -1) | Out-Null\n$targetRange.Text = \"Version management help many developers to collaborate while working on same project.\"\n\n# Apply bold + red formatting to the whole paragraph (run text AND the\n# paragraph mark) so both the run and the paragraph's own rPr pick up\n# the formatting.\n$wdColorRed = [Microsoft.Office.Interop.Word.WdColor]::wdColorRed\n$fullParagraphRange = $targetParagraph.Range\n$fullParagraphRange.Font.Bold = $true\n$fullParagraphRange.Font.Color = $wdColorRed\n\n# Re-create the \"_GoBack\" bookmark at the very last paragraph in the\n# document (a collapsed bookmark at that paragraph's range).\n$paragraphCount = $d.Paragraphs.Count\n$lastParagraph = $d.Paragraphs($paragraphCount)\n$d.Bookmarks.Add(\"_GoBack\", $lastParagraph.Range) | Out-Null\n"}
